{"js": "// The diff replaces the text of 25 table cells (5 populated rows x 5\n// columns) in the single table of the document with new \"a\u00f7b=c, d\"\n// strings, leaving everything else (the date header paragraph and the\n// empty spacer rows) untouched.\n//\n// Each cell is addressed positionally (row, column) rather than by\n// matching the old text, because several of the old values repeat\n// (e.g. \"83\u00f77=11, 6\" occurs twice) but map to different replacements\n// depending on their position.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newText] in document order.\nconst replacements = [\n  [0, 0, \"14\u00f78=1, 6\"],\n  [0, 1, \"63\u00f73=21, 0\"],\n  [0, 2, \"92\u00f78=11, 4\"],\n  [0, 3, \"50\u00f74=12, 2\"],\n  [0, 4, \"32\u00f72=16, 0\"],\n\n  [4, 0, \"59\u00f77=8, 3\"],\n  [4, 1, \"48\u00f79=5, 3\"],\n  [4, 2, \"65\u00f76=10, 5\"],\n  [4, 3, \"41\u00f77=5, 6\"],\n  [4, 4, \"95\u00f79=10, 5\"],\n\n  [8, 0, \"64\u00f73=21, 1\"],\n  [8, 1, \"46\u00f75=9, 1\"],\n  [8, 2, \"73\u00f74=18, 1\"],\n  [8, 3, \"84\u00f75=16, 4\"],\n  [8, 4, \"84\u00f75=16, 4\"],\n\n  [12, 0, \"21\u00f76=3, 3\"],\n  [12, 1, \"36\u00f72=18, 0\"],\n  [12, 2, \"47\u00f78=5, 7\"],\n  [12, 3, \"13\u00f78=1, 5\"],\n  [12, 4, \"29\u00f76=4, 5\"],\n\n  [16, 0, \"52\u00f75=10, 2\"],\n  [16, 1, \"82\u00f73=27, 1\"],\n  [16, 2, \"17\u00f76=2, 5\"],\n  [16, 3, \"69\u00f79=7, 6\"],\n  [16, 4, \"70\u00f77=10, 0\"],\n];\n\nfor (const [row, col, text] of replacements) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# The diff replaces the text of 25 table cells (5 populated rows x 5\n# columns) in the single table of the document with new \"a\u00f7b=c, d\"\n# strings, leaving everything else (the date header paragraph and the\n# empty spacer rows) untouched.\n#\n# Each cell is addressed positionally (row, column) rather than by\n# matching the old text, because several of the old values repeat\n# (e.g. \"83\u00f77=11, 6\" occurs twice) but map to different replacements\n# depending on their position.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# row, column (1-based, matching Word COM's Table.Cell indexing), new text\n$replacements = @(\n    @(1, 1, \"14\u00f78=1, 6\"),\n    @(1, 2, \"63\u00f73=21, 0\"),\n    @(1, 3, \"92\u00f78=11, 4\"),\n    @(1, 4, \"50\u00f74=12, 2\"),\n    @(1, 5, \"32\u00f72=16, 0\"),\n\n    @(5, 1, \"59\u00f77=8, 3\"),\n    @(5, 2, \"48\u00f79=5, 3\"),\n    @(5, 3, \"65\u00f76=10, 5\"),\n    @(5, 4, \"41\u00f77=5, 6\"),\n    @(5, 5, \"95\u00f79=10, 5\"),\n\n    @(9, 1, \"64\u00f73=21, 1\"),\n    @(9, 2, \"46\u00f75=9, 1\"),\n    @(9, 3, \"73\u00f74=18, 1\"),\n    @(9, 4, \"84\u00f75=16, 4\"),\n    @(9, 5, \"84\u00f75=16, 4\"),\n\n    @(13, 1, \"21\u00f76=3, 3\"),\n    @(13, 2, \"36\u00f72=18, 0\"),\n    @(13, 3, \"47\u00f78=5, 7\"),\n    @(13, 4, \"13\u00f78=1, 5\"),\n    @(13, 5, \"29\u00f76=4, 5\"),\n\n    @(17, 1, \"52\u00f75=10, 2\"),\n    @(17, 2, \"82\u00f73=27, 1\"),\n    @(17, 3, \"17\u00f76=2, 5\"),\n    @(17, 4, \"69\u00f79=7, 6\"),\n    @(17, 5, \"70\u00f77=10, 0\")\n)\n\nforeach ($item in $replacements) {\n    $row = $item[0]\n    $col = $item[1]\n    $text = $item[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
